$wb = $excel.ActiveWorkbook

$oldGuid = "03e11c4c-a08f-4387-9923-82bd72a53dbf"
$newGuid = "e547055f-1636-439f-b445-4340d066b4e5"

$oldHoHash = "dc0d75ec29b062cf252d2f097c0ae1fa7f8445ca"
$newHoHash = "3f52c374c44e51373916dff20ed5fc3b032238d5"

# ---------------------------------------------------------------
# Overview sheet: B2 hyperlink display text picks up the new guid.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"

# Keep the existing hyperlink target (rId2 unchanged) but refresh the
# display text shown in the cell/link tooltip.
$ovHl = $wsOverview.Hyperlinks.Item(1)
$ovHl.TextToDisplay = "e2e\" + $newGuid + ".md"

# "Latest HO Xliff Generate Date" (also shared with de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-13 07:15:47"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Source File Name (A2) + its hyperlink display text
$wsZh.Range("A2").Value = $newGuid + ".md"

# Drop the "Latest Target File" hyperlink (I2) entirely - the new report
# no longer links/populates this column.
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/4f0edb507f4ad80c18368c002534df682c74d621/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = $wsZh.Range("J2").Style

$wsZh.Range("J2").Value = ""

$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newGuid + ".md"

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/4f0edb507f4ad80c18368c002534df682c74d621/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = $wsDe.Range("J2").Style

$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------
# Shared "Latest Handoff File" / generated xlf names move to the new guid
# ---------------------------------------------------------------
$wsZh.Range("G2").Value = $newGuid + "." + $newHoHash + ".zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-13 07:15:40"

$wsDe.Range("G2").Value = $newGuid + "." + $newHoHash + ".de-de.xlf"
